$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E (old E "Mã GV" / DSSV-index column shifts to F),
# then fill the new column E with the "Mã GVU" header + GVU#### values.
$ws.Columns("E:E").Insert()

$ws.Range("E1").Value = "Mã GVU"
$ws.Range("E2").Value = "GVU1234"
$ws.Range("E3").Value = "GVU1235"
$ws.Range("E4").Value = "GVU1236"
$ws.Range("E5").Value = "GVU1237"

# Match the column widths captured by Excel's own best-fit sizing
# (target raw widths: B=11.21875, E=8.6640625, F=36.33203125; the inputs
# below are tuned so this engine's width rounding lands as close as
# possible to those stored values).
$ws.Columns("B:B").ColumnWidth = 10.2538
$ws.Columns("E:E").ColumnWidth = 7.7541
$ws.Columns("F:F").ColumnWidth = 35.417

# Selection now spans the whole used range instead of the old E4:E5 pick.
$ws.Range("A1:F5").Select()

# Page setup was touched (portrait orientation explicitly stamped).
$ws.PageSetup.Orientation = 1

Write-Host "edit applied"
